$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellAddress, $TextValue)
    $r = $ws.Range($CellAddress)
    $r.NumberFormat = "@"
    $r.Value = $TextValue
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "61.194.15"
$ws.Range("E2").Value = "  +5.72%  "
$ws.Range("D3").Value = "2.385.28"
$ws.Range("E3").Value = "  +4.30%  "
Set-TextValue "D4" "0.997"
$ws.Range("E4").Value = "  -0.25%  "
Set-TextValue "D5" "552.00"
$ws.Range("E5").Value = "  +2.41%  "
Set-TextValue "D6" "135.18"
$ws.Range("E6").Value = "  +3.18%  "
Set-TextValue "D7" "0.998"
$ws.Range("E7").Value = "  -0.26%  "
Set-TextValue "D8" "0.593"
$ws.Range("E8").Value = "  +4.62%  "
$ws.Range("D9").Value = "2.380.83"
$ws.Range("E9").Value = "  +4.20%  "
Set-TextValue "D10" "0.102"
$ws.Range("E10").Value = "  +1.88%  "
Set-TextValue "D11" "5.57"
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("E12").Value = "  +1.68%  "
Set-TextValue "D13" "0.341"
$ws.Range("E13").Value = "  +3.11%  "
Set-TextValue "D14" "24.41"
$ws.Range("E14").Value = "  +4.04%  "
$ws.Range("D15").Value = "2.787.58"
$ws.Range("E15").Value = "  +3.51%  "
$ws.Range("D16").Value = "61.021.05"
$ws.Range("E16").Value = "  +5.44%  "
Set-TextValue "D17" "0.0000135"
$ws.Range("E17").Value = "  +2.69%  "
$ws.Range("D18").Value = "2.339.58"
$ws.Range("E18").Value = "  +3.54%  "
Set-TextValue "D19" "10.87"
$ws.Range("E19").Value = "  +2.68%  "
Set-TextValue "D20" "4.27"
$ws.Range("E20").Value = "  +0.83%  "
Set-TextValue "D21" "6.91"
$ws.Range("E21").Value = "  +7.93%  "
Set-TextValue "D22" "321.06"
$ws.Range("E22").Value = "  +2.58%  "
$ws.Range("E23").Value = "  +0.67%  "
Set-TextValue "D24" "63.62"
$ws.Range("E24").Value = "  +1.01%  "
Set-TextValue "D25" "0.175"
$ws.Range("E25").Value = "  +5.08%  "
Set-TextValue "D26" "0.996"
$ws.Range("E26").Value = "  -0.42%  "
Set-TextValue "D27" "8.23"
$ws.Range("E27").Value = "  +3.23%  "
Set-TextValue "D28" "1.36"
$ws.Range("E28").Value = "  +5.29%  "
Set-TextValue "D29" "1.76"
$ws.Range("E29").Value = "  +3.07%  "
$ws.Range("D30").Value = "0.0₃0761"
$ws.Range("E30").Value = "  +5.73%  "
Set-TextValue "D31" "171.90"
$ws.Range("E31").Value = "  +1.27%  "
Set-TextValue "D32" "1.15"
$ws.Range("E32").Value = "  +6.89%  "
Set-TextValue "D33" "6.00"
$ws.Range("E33").Value = "  +4.97%  "
Set-TextValue "D34" "1.44"
$ws.Range("E34").Value = "  +16.35%  "
Set-TextValue "D35" "0.390"
$ws.Range("E35").Value = "  +2.79%  "
Set-TextValue "D36" "18.18"
$ws.Range("E36").Value = "  +2.89%  "
Set-TextValue "D37" "0.998"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D38" "4.22"
$ws.Range("E38").Value = "  +8.37%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D39" "0.999"
$ws.Range("E39").Value = "  -0.08%  "
Set-TextValue "D40" "328.51"
$ws.Range("E40").Value = "  +13.31%  "
Set-TextValue "D41" "1.59"
$ws.Range("E41").Value = "  +7.18%  "
Set-TextValue "D42" "38.52"
$ws.Range("E42").Value = "  +1.85%  "
Set-TextValue "D43" "146.37"
$ws.Range("E43").Value = "  +3.80%  "
Set-TextValue "D44" "3.52"
$ws.Range("E44").Value = "  +3.92%  "
Set-TextValue "D45" "0.0962"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D46" "19.58"
$ws.Range("E46").Value = "  +7.59%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D47" "0.0506"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D48" "0.569"
$ws.Range("E48").Value = "  +3.12%  "
Set-TextValue "D49" "0.0216"
$ws.Range("E49").Value = "  +2.23%  "
Set-TextValue "D50" "11.02"
$ws.Range("E50").Value = "  +0.58%  "
Set-TextValue "D51" "1.57"
$ws.Range("E51").Value = "  +5.89%  "

Write-Host "Applied 108 cell updates"
